# Update "paises" (countries) COVID dashboard with newer figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 18:39"

# Two countries in the ranked list swapped places because their totals
# changed relative to their neighbour. Swap the country names so that the
# row with the higher (updated) total keeps sitting above the other one.
$ws.Range("A33").Value = "Ecuador"
$ws.Range("A34").Value = "Marruecos"

$ws.Range("A40").Value = "Chequia"
$ws.Range("A41").Value = "Kuwait"

$ws.Range("A135").Value = "Sri Lanka"
$ws.Range("A136").Value = "Siria"

# Refresh the numeric figures (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Muertes hoy, Muertes) for the rows whose data changed.

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7912262
$ws.Range("C4").Value = 17784
$ws.Range("D4").Value = 5068912
$ws.Range("E4").Value = 2624472
$ws.Range("G4").Value = 230
$ws.Range("H4").Value = 218878

# Row 23 - Turquia
$ws.Range("B23").Value = 334031
$ws.Range("C23").Value = 1649
$ws.Range("D23").Value = 293145
$ws.Range("E23").Value = 32108
$ws.Range("G23").Value = 56
$ws.Range("H23").Value = 8778

# Row 29 - Canada
$ws.Range("B29").Value = 180026
$ws.Range("C29").Value = 1909
$ws.Range("D29").Value = 151335
$ws.Range("E29").Value = 19085
$ws.Range("G29").Value = 21
$ws.Range("H29").Value = 9606

# Row 33 - Ecuador (now ranked above Marruecos)
$ws.Range("B33").Value = 146828
$ws.Range("C33").Value = 980
$ws.Range("D33").Value = 120511
$ws.Range("E33").Value = 14129
$ws.Range("G33").Value = 13
$ws.Range("H33").Value = 12188

# Row 34 - Marruecos (now ranked below Ecuador)
$ws.Range("B34").Value = 146398
$ws.Range("D34").Value = 123022
$ws.Range("E34").Value = 20846
$ws.Range("H34").Value = 2530

# Row 40 - Chequia (now ranked above Kuwait)
$ws.Range("B40").Value = 111978
$ws.Range("C40").Value = 2604
$ws.Range("D40").Value = 53202
$ws.Range("E40").Value = 57842
$ws.Range("G40").Value = 29
$ws.Range("H40").Value = 934

# Row 41 - Kuwait (now ranked below Chequia)
$ws.Range("B41").Value = 110568
$ws.Range("C41").Value = 492
$ws.Range("D41").Value = 102722
$ws.Range("E41").Value = 7191
$ws.Range("G41").Value = 6
$ws.Range("H41").Value = 655

# Row 62 - Singapur
$ws.Range("D62").Value = 57698
$ws.Range("E62").Value = 141

# Row 87 - Grecia
$ws.Range("B87").Value = 22078
$ws.Range("C87").Value = 306
$ws.Range("E87").Value = 11653
$ws.Range("G87").Value = 5
$ws.Range("H87").Value = 436

# Row 100 - Montenegro
$ws.Range("B100").Value = 13641
$ws.Range("C100").Value = 293
$ws.Range("D100").Value = 9429
$ws.Range("E100").Value = 4014
$ws.Range("G100").Value = 4
$ws.Range("H100").Value = 198

# Row 104 - Guinea
$ws.Range("B104").Value = 10996
$ws.Range("C104").Value = 42
$ws.Range("D104").Value = 10304
$ws.Range("E104").Value = 623
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 69

# Row 110 - Uganda
$ws.Range("B110").Value = 9701
$ws.Range("C110").Value = 163
$ws.Range("D110").Value = 6022
$ws.Range("E110").Value = 3586
$ws.Range("G110").Value = 7
$ws.Range("H110").Value = 93

# Row 123 - Malaui
$ws.Range("B123").Value = 5821
$ws.Range("C123").Value = 8
$ws.Range("D123").Value = 4644
$ws.Range("E123").Value = 997

# Row 135 - Sri Lanka (now ranked above Siria)
$ws.Range("B135").Value = 4628
$ws.Range("C135").Value = 105
$ws.Range("D135").Value = 3306
$ws.Range("E135").Value = 1309
$ws.Range("H135").Value = 13

# Row 136 - Siria (now ranked below Sri Lanka)
$ws.Range("B136").Value = 4616
$ws.Range("D136").Value = 1235
$ws.Range("E136").Value = 3163
$ws.Range("H136").Value = 218

# Row 156 - Sierra Leona
$ws.Range("B156").Value = 2300
$ws.Range("C156").Value = 5
$ws.Range("D156").Value = 1725
$ws.Range("E156").Value = 503

# Row 160 - Republica de Chipre
$ws.Range("B160").Value = 1984
$ws.Range("C160").Value = 33
$ws.Range("E160").Value = 590

# Row 177 - Burundi
$ws.Range("B177").Value = 517
$ws.Range("C177").Value = 2
$ws.Range("E177").Value = 44
